$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the contents of row 7 (A7:I7) but keep existing cell formatting,
# so the "free air space" printed in row 7 goes away (fixing printing
# free air space), matching the pattern of the already-blank rows below.
$ws.Range("A7:I7").ClearContents()
